# Dataset_Rotas_BI.xlsx edit
# - Clear stray "25" placeholder values out of column J (rows 2-247, except the
#   four rows that receive real percentage data below).
# - Move the vehicle info that used to live in F248/G248 and in the three
#   "orphan" rows 249-251 up into F237:J240 (their proper rows), which also
#   carries the right number formats (copied from the original rows 249-251).
# - Remove the now empty trailing rows 249-251.
# - Clear F248/G248/J248 (their data moved up to row 237).
# - Update the active selection to match the final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 247

# 1) Clear the "25" filler values from column J for the regular data rows,
#    but skip rows 238-240 (they get real percentage data below) since those
#    are handled explicitly afterwards. Row 237 also loses its "25" (it does
#    not receive a replacement value).
for ($r = 2; $r -le $lastDataRow; $r++) {
    if ($r -ge 238 -and $r -le 240) { continue }
    $ws.Cells.Item($r, 10).ClearContents()
}

# 2) Copy the number formatting for F/G/J from the source rows (249-251) onto
#    the destination rows (238-240) before we touch their values, so the
#    styles (borders/number formats) match what the diff expects.
$ws.Cells.Item(249, 6).Copy()
$ws.Cells.Item(238, 6).PasteSpecial(-4122)
$ws.Cells.Item(249, 7).Copy()
$ws.Cells.Item(238, 7).PasteSpecial(-4122)
$ws.Cells.Item(249, 10).Copy()
$ws.Cells.Item(238, 10).PasteSpecial(-4122)

$ws.Cells.Item(250, 6).Copy()
$ws.Cells.Item(239, 6).PasteSpecial(-4122)
$ws.Cells.Item(250, 7).Copy()
$ws.Cells.Item(239, 7).PasteSpecial(-4122)
$ws.Cells.Item(250, 10).Copy()
$ws.Cells.Item(239, 10).PasteSpecial(-4122)

$ws.Cells.Item(251, 6).Copy()
$ws.Cells.Item(240, 6).PasteSpecial(-4122)
$ws.Cells.Item(251, 7).Copy()
$ws.Cells.Item(240, 7).PasteSpecial(-4122)
$ws.Cells.Item(251, 10).Copy()
$ws.Cells.Item(240, 10).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Fill in row 237 with the vehicle data that used to sit in F248/G248
#    (keeps its original style, nothing to copy there).
$ws.Cells.Item(237, 6).Value = "Caminhão Toco"
$ws.Cells.Item(237, 7).Value = 4

# 4) Fill in rows 238-240 with the vehicle data that used to sit in the
#    orphan rows 249-251.
$ws.Cells.Item(238, 6).Value = "Carro 1.0"
$ws.Cells.Item(238, 7).Value = 13
$ws.Cells.Item(238, 10).Value = 0.9

$ws.Cells.Item(239, 6).Value = "Sprinter Van 516 Standard"
$ws.Cells.Item(239, 7).Value = 11
$ws.Cells.Item(239, 10).Value = 9

$ws.Cells.Item(240, 6).Value = "Moto 110–125"
$ws.Cells.Item(240, 7).Value = 40
$ws.Cells.Item(240, 10).Value = 0.15

# 5) The data that used to live in F248/G248 has moved to row 237, so clear
#    it (along with its own "25" filler in J248).
$ws.Cells.Item(248, 6).ClearContents()
$ws.Cells.Item(248, 7).ClearContents()
$ws.Cells.Item(248, 10).ClearContents()

# 6) Remove the now-obsolete orphan rows (shifts dimension down to O248).
$ws.Rows("249:251").Delete()

# 7) Update the current selection to match the final workbook state.
$ws.Activate()
$ws.Range("H246").Select() | Out-Null
